$d = $word.ActiveDocument

# Builds a WordOpenXML "flat OPC" fragment containing a single paragraph
# whose text is split into one run per word plus one run per inter-word
# space (matching the target structure), and inserts it in place of the
# supplied paragraph's range. Using InsertXML (raw OOXML) rather than
# simple text assignment is what lets the individually-styled-identically
# runs remain distinct <w:r> elements instead of being recombined into a
# single run when the document is serialized.
function Set-ParagraphWordRuns($paragraph, $pStyle, $text) {
    $words = $text.Split(" ")
    $runsXml = ""
    for ($i = 0; $i -lt $words.Length; $i++) {
        if ($i -gt 0) {
            $runsXml += '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
        }
        $runsXml += '<w:r><w:t xml:space="preserve">' + $words[$i] + '</w:t></w:r>'
    }

    $pPrXml = ""
    if ($pStyle) {
        $pPrXml = '<w:pPr><w:pStyle w:val="' + $pStyle + '"/></w:pPr>'
    }

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $pPrXml + $runsXml + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $null = $paragraph.Range.InsertXML($xml)
}

function Get-ParagraphByStyle($doc, $styleName) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Style.NameLocal -eq $styleName) {
            return $p
        }
    }
    return $null
}

$titlePara = Get-ParagraphByStyle $d "Title"
Set-ParagraphWordRuns $titlePara "Title" "Questions: Introduction to radians"

$authorPara = Get-ParagraphByStyle $d "Author"
Set-ParagraphWordRuns $authorPara "Author" "Mark Toner, Ifan Howell-Baines"

$abstractPara = Get-ParagraphByStyle $d "Abstract"
Set-ParagraphWordRuns $abstractPara "Abstract" "Questions relating to the introduction to radians study guide."
